$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: D20 1 -> 2
$ws.Range("D20").Value2 = 2

# New dialog rows (22-27) - Fleur character added
# Column G (sprite path) set first across all affected rows
$ws.Range("G22:G27").Value2 = "Sprite/TalkSprite/Fleur"

# Column E (speaker name) for rows 22, 24-27 (row 23 uses {Name} placeholder)
$ws.Range("E22").Value2 = "플뢰르"
$ws.Range("E24").Value2 = "플뢰르"
$ws.Range("E25").Value2 = "플뢰르"
$ws.Range("E26").Value2 = "플뢰르"
$ws.Range("E27").Value2 = "플뢰르"
$ws.Range("E23").Value2 = "{Name}"

# Column F (dialog text)
$ws.Range("F23").Value2 = "미안.. 갑자기 데자뷰가.."
$ws.Range("F22").Value2 = "뭐..뭔가? 갑자기…"
$ws.Range("F24").Value2 = "겨우 저 정도 적들에게 쓰러지다니..`n내 체면이 말이 아니군…"
$ws.Range("F25").Value2 = "본인의 능력치 강화나 스크롤은 제대로 구한건가?"
$ws.Range("F26").Value2 = "우두머리들은 강력하니 그때 만큼은`n직접 몸을 조종하는게 좋을거야."
$ws.Range("F27").Value2 = "그럼 다시 한 번 가보자고."

# Column J (text color)
$ws.Range("J22").Value2 = "#FFFFFF"
$ws.Range("J23").Value2 = "#606060"
$ws.Range("J24").Value2 = "#FFFFFF"
$ws.Range("J25").Value2 = "#FFFFFF"
$ws.Range("J26").Value2 = "#FFFFFF"
$ws.Range("J27").Value2 = "#FFFFFF"

# Wrap text + row height for multi-line dialog rows
$ws.Range("F24").WrapText = $true
$ws.Range("F26").WrapText = $true
$ws.Rows.Item(24).RowHeight = 33
$ws.Rows.Item(26).RowHeight = 33

# Rows 28-30: chapter/cycle index change (B: 2->3, C renumbered 11,12,13 -> 1,2,3)
$ws.Range("B28").Value2 = 3
$ws.Range("C28").Value2 = 1
$ws.Range("B29").Value2 = 3
$ws.Range("C29").Value2 = 2
$ws.Range("B30").Value2 = 3
$ws.Range("C30").Value2 = 3

# Selection change
$ws.Range("F23").Select()
